$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.266.51"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "3.582.69"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.23%  "
$ws.Range("D7").Value = "3.577.03"
$ws.Range("E7").Value = "  +0.61%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.412"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "4.192.34"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000205"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("D16").Value = "3.595.39"
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("D18").Value = "66.374.49"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "423.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.608"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.71%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000120"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "3.582.08"
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("E32").Value = "  +3.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "24.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("E34").Value = "  -1.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.72"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("E38").Value = "  -3.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "174.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0850"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.879"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "45.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.54%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.946"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.90%  "
